$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "department" column (C) from the single faculty name
# "FACULTY OF BUSINESS & TECHNOLOGY" to more specific category names.
$ws.Range("C2").Value = "Business"
$ws.Range("C3").Value = "Business"
$ws.Range("C4").Value = "Business"
$ws.Range("C5").Value = "Business"
$ws.Range("C6").Value = "Business"
$ws.Range("C7").Value = "Business"
$ws.Range("C8").Value = "Business"
$ws.Range("C9").Value = "Business"
$ws.Range("C10").Value = "Business"
$ws.Range("C11").Value = "Business"
$ws.Range("C12").Value = "Business"

$ws.Range("C13").Value = "Information Technology"
$ws.Range("C14").Value = "Information Technology"
$ws.Range("C15").Value = "Information Technology"

$ws.Range("C16").Value = "Building and Construction"

$ws.Range("C17").Value = "Packages"
$ws.Range("C18").Value = "Packages"
$ws.Range("C19").Value = "Packages"
$ws.Range("C20").Value = "Packages"
$ws.Range("C21").Value = "Packages"
$ws.Range("C22").Value = "Packages"
